$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    # Force the cell to keep a text (shared-string) type even though the
    # string looks like a plain number (e.g. "-0.041"), matching how these
    # report tables store numeric-looking results as text. Temporarily mark
    # the cell as Text, assign, then restore the default ("Normal") style so
    # no visible formatting change is left behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Column headers / row labels: drop the " Diff-in-Diff" suffix
$ws.Range("A2").Value = "FFR"
$ws.Range("A3").Value = "C/A"
$ws.Range("A4").Value = "U"
$ws.Range("A5").Value = '$\pi$'

$ws.Range("B1").Value = "FFR"
$ws.Range("C1").Value = "C/A"
$ws.Range("D1").Value = "U"
$ws.Range("E1").Value = '$\pi$'

# Updated coefficient values (text cells containing "***" already stay text
# automatically since they aren't parsed as plain numbers)
Set-TextValue $ws.Range("B4") "-0.041"
Set-TextValue $ws.Range("B5") "0.036"
Set-TextValue $ws.Range("B6") "-0.135"

$ws.Range("C2").Value = "3.116***"
$ws.Range("C4").Value = "-5.92***"
$ws.Range("C5").Value = "6.726***"
Set-TextValue $ws.Range("C6") "1.098"

Set-TextValue $ws.Range("D2") "-0.029"
$ws.Range("D3").Value = "-0.043***"
Set-TextValue $ws.Range("D5") "-0.086"
Set-TextValue $ws.Range("D6") "0.023"

Set-TextValue $ws.Range("E2") "0.014"
$ws.Range("E3").Value = "0.027***"
Set-TextValue $ws.Range("E4") "-0.048"
Set-TextValue $ws.Range("E6") "-0.045"
